$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

# E19 was stored as text "533758"; convert it to a real number.
$ws.Range("E19").Value = 533758

# Append a new row 20 with the same stock repeated, at a later timestamp.
$ws.Range("A20").Value = "13/06/2024 07:45:46"
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = "APLAPOLLO"
$ws.Range("D20").Value = "Apl Apollo Tubes Limited"
# bsecode for row 20 stays a text value (as row 19 originally was), so force text with a leading apostrophe.
$ws.Range("E20").Value = "'533758"
$ws.Range("F20").Value = -2.33
$ws.Range("G20").Value = 1544.05
$ws.Range("H20").Value = 363103
